$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Chapter 2 (column C) index-term checkmarks for rows 23-46 ("check"/"n/a"),
# mirroring the pattern already present for rows 3-22.
$chapter2 = [ordered]@{
    23 = "check"
    24 = "check"
    25 = "n/a"
    26 = "check"
    27 = "n/a"
    28 = "n/a"
    29 = "n/a"
    30 = "n/a"
    31 = "n/a"
    32 = "n/a"
    33 = "n/a"
    34 = "n/a"
    35 = "n/a"
    36 = "n/a"
    37 = "n/a"
    38 = "n/a"
    39 = "check"
    40 = "check"
    41 = "n/a"
    42 = "n/a"
    43 = "n/a"
    44 = "n/a"
    45 = "check"
    46 = "check"
}

foreach ($row in $chapter2.Keys) {
    $ws.Cells.Item($row, 3).Value = $chapter2[$row]
}

# Match the author's final cursor position/selection.
$ws.Range("C46").Select()
